# Biolog Daten PM2A Platten.xlsx - "update umay substrate grwoth"
#
# The real edit behind this commit is on the last worksheet ("31-05-2021"):
# a block of 6 blank rows that separated the first measurement table
# (rows 2:10) from the second one (rows 28:36) was deleted, so the second
# table moves up to rows 22:30 - matching the row layout already used on
# the "21-05-2021" / "26-05-2021" sheets (2:10, 12:20, 22:30).
#
# The rest of the diff (selection/active-cell bookkeeping, the active tab)
# is just Excel recording where the user's cursor ended up while doing
# that edit, so we reproduce that too.

$wb = $excel.ActiveWorkbook

$wsOverview  = $wb.Worksheets.Item("Overview")
$ws2105      = $wb.Worksheets.Item("21-05-2021")
$ws2605      = $wb.Worksheets.Item("26-05-2021")
$ws3105      = $wb.Worksheets.Item("31-05-2021")

# --- The actual data edit -------------------------------------------------
# Remove the 6 superfluous blank rows between the two data tables on the
# "31-05-2021" sheet (rows 11:16 are blank; deleting them shifts the lower
# table from 28:36 up to 22:30, same as rows 28:36 do not need to stay in
# place since nothing else references them).
$ws3105.Rows("11:16").Delete()

# --- View / selection state -------------------------------------------
# Sheet "21-05-2021": cursor ends up on N14 (scrolled so row 7 is on top).
$ws2105.Activate()
$ws2105.Range("N14").Select()
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 7
    $win.ScrollColumn = 1
} catch {}

# Sheet "26-05-2021": no change to its stored selection (still P29).
$ws2605.Activate()
$ws2605.Range("P29").Select()

# Sheet "Overview": no change to its stored selection (still H15).
$wsOverview.Activate()
$wsOverview.Range("H15").Select()

# Sheet "31-05-2021" ends up as the active tab, cursor on H35 (scrolled so
# row 14 is on top) - this becomes the workbook's active sheet on save.
$ws3105.Activate()
$ws3105.Range("H35").Select()
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 14
    $win.ScrollColumn = 1
} catch {}
